# Weekly update: insert two new price-report rows (week of 2023-12-15,
# serial 45275) for "Apio" at "Terminal Hortofrutícola Agro Chillán",
# pushing the existing rows 422-510 down to 424-512.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 422 (old data shifts down).
$ws.Rows.Item(422).EntireRow.Insert()
$ws.Rows.Item(423).EntireRow.Insert()

# New row 422 - "Primera" quality
$ws.Range("A422").Value = 7
$ws.Range("B422").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C422").Value = "Ñuble"
$ws.Range("D422").Value = 45275
$ws.Range("E422").Value = 16
$ws.Range("F422").Value = 100112017
$ws.Range("G422").Value = "Apio"
$ws.Range("H422").Value = "Americana (o)"
$ws.Range("I422").Value = "Primera"
$ws.Range("J422").Value = 100
$ws.Range("K422").Value = 12000
$ws.Range("L422").Value = 12000
$ws.Range("M422").Value = 12000
$ws.Range("N422").Value = "$/docena de matas"
$ws.Range("O422").Value = "Provincia del Elquí"
$ws.Range("P422").Value = 2000
$ws.Range("Q422").Value = 6
$ws.Range("R422").Value = "Hortaliza"

# New row 423 - "Segunda" quality
$ws.Range("A423").Value = 7
$ws.Range("B423").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C423").Value = "Ñuble"
$ws.Range("D423").Value = 45275
$ws.Range("E423").Value = 16
$ws.Range("F423").Value = 100112017
$ws.Range("G423").Value = "Apio"
$ws.Range("H423").Value = "Americana (o)"
$ws.Range("I423").Value = "Segunda"
$ws.Range("J423").Value = 100
$ws.Range("K423").Value = 9000
$ws.Range("L423").Value = 9000
$ws.Range("M423").Value = 9000
$ws.Range("N423").Value = "$/docena de matas"
$ws.Range("O423").Value = "Provincia del Elquí"
$ws.Range("P423").Value = 1500
$ws.Range("Q423").Value = 6
$ws.Range("R423").Value = "Hortaliza"
